$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data per upstream scrape.
# Values that look like plain numbers must be forced to Text so Excel
# doesn't silently coerce them (and drop formatting such as trailing zeros
# or the thousands-dot notation used by this sheet), matching how the
# source data is stored (inline text, not numeric).

$ws.Range('D2').Value = '70.839.99'
$ws.Range('E2').Value = '  -2.02%  '
$ws.Range('D3').Value = '3.638.99'
$ws.Range('E3').Value = '  +0.47%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '584.32'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -2.35%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '176.29'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.54%  '
$ws.Range('D7').Value = '3.635.19'
$ws.Range('E7').Value = '  +0.65%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.615'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('E9').Value = '  -0.07%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.197'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -4.82%  '
$ws.Range('E11').Value = '  +16.06%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.608'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.35%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '48.61'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -4.23%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.0000285'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.94%  '
$ws.Range('D15').Value = '4.224.14'
$ws.Range('E15').Value = '  +0.47%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '676.09'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -4.23%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '8.96'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').Value = '3.630.07'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').Value = '70.839.13'
$ws.Range('E19').Value = '  -2.16%  '
$ws.Range('E20').Value = '  -0.40%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '17.82'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -4.48%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '11.52'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.99%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.941'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.64%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '17.19'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -4.08%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '100.12'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -5.06%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '3.93'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -2.81%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '2.80'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -2.17%  '
$ws.Range('E28').Value = '  +0.03%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.85'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.62%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '34.70'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -2.50%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '9.15'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.18%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.31'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -5.41%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '7.59'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.57%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.39'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -6.67%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '4.00'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -4.69%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '576.44'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -3.18%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '11.11'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -2.36%  '
$ws.Range('E38').Value = '  -0.94%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '58.57'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('E40').Value = '  -0.06%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0455'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').Value = '3.569.03'
$ws.Range('E42').Value = '  -1.99%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.346'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.97%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.140'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -3.03%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '34.50'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -4.10%  '
$ws.Range('D46').Value = '0.0₃0734'
$ws.Range('E46').Value = '  -5.90%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.69'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -4.06%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.90'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +3.12%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.134'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.78%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '137.37'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +2.50%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.91'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -3.12%  '
